# BrightCoffee Presentation edit script
# Applies:
#   1. Slide 5 (sldId 274) title "Sales VS Revenue per Product Type":
#        - move textbox up (y offset 292625 -> 211325)
#        - bold the title run and set font to Fira Sans Extra Condensed
#   2. Slide 6 (sldId 277) title "Revenue per time of the day":
#        - bold the title run and set font to Fira Sans Extra Condensed
#   3. Slide 7 (sldId 282) title box (was "COFFEE INFOGRAPHICS"):
#        - add a new first line "Sales per store location" (bold, centered,
#          Fira Sans Extra Condensed)
#        - remove the old "COFFEE INFOGRAPHICS" wording
#   4. Slide 8 (sldId 288) title "Number of Sales vs day of the week":
#        - bold the title run and set font to Fira Sans Extra Condensed
#        - add a new empty textbox (leftover placeholder) after the last shape

$p = $ppt.ActivePresentation

$fontName = "Fira Sans Extra Condensed"

# Shape.Top/.Left/.Width/.Height are expressed in points (1 pt = 12700 EMU)
# through this COM surface, same as real PowerPoint VBA automation.
function EmuToPt($emu) { return $emu / 12700.0 }

# ---------------------------------------------------------------------------
# 1. Slide 5 - "Sales VS Revenue per Product Type"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$sh5.Top = EmuToPt 211325
$run5 = $sh5.TextFrame.TextRange.Paragraphs(1, 1)
$run5.Font.Bold = $true
$run5.Font.Name = $fontName

# ---------------------------------------------------------------------------
# 2. Slide 6 - "Revenue per time of the day"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(1)
$run6 = $sh6.TextFrame.TextRange.Paragraphs(1, 1)
$run6.Font.Bold = $true
$run6.Font.Name = $fontName

# ---------------------------------------------------------------------------
# 3. Slide 7 - "COFFEE INFOGRAPHICS" -> "Sales per store location"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(1)
$tr7 = $sh7.TextFrame.TextRange
$tr7.InsertBefore("Sales per store location`r")
$newPara7 = $tr7.Paragraphs(1, 1)
$newPara7.ParagraphFormat.Alignment = 2
$newPara7.Font.Bold = $true
$newPara7.Font.Name = $fontName
$newPara7.Font.Size = 28

$oldPara7 = $tr7.Paragraphs(2, 1)
$oldPara7.Text = "`r"

# ---------------------------------------------------------------------------
# 4. Slide 8 - "Number of Sales vs day of the week"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(1)
$run8 = $sh8.TextFrame.TextRange.Paragraphs(1, 1)
$run8.Font.Bold = $true
$run8.Font.Name = $fontName

# New empty textbox added at the end of slide 8's shape tree
$newBox = $s8.Shapes.AddTextbox(1, (EmuToPt 554792), (EmuToPt 777282), (EmuToPt 8520600), (EmuToPt 572700))
$newBox.Name = "Google Shape;2587;p40"
$newBox.TextFrame.WordWrap = $true
$newBox.Fill.Visible = $false
$newBox.Line.Visible = $false
